# Adds role/assignment descriptions to each team member on the "Clanovi tima"
# slide (slide 3), fixes a couple of misspelled names, as per:
#   "Dodane uloge u prezentaciju"
#
# Character offsets below are computed against the ORIGINAL text of the
# content placeholder:
#   "Filip Ljubotina\rMarko Pavić\rMihael Breznicki-Herceg\rLara Ćorić\r
#    Ana Vuksanović\rKatarina Klaric\rNoa Milin"
# Para 1 "Filip Ljubotina"           start=1  len=15
# Para 2 "Marko Pavić"               start=17 len=11
# Para 3 "Mihael Breznicki-Herceg"   start=29 len=23  ("Breznicki-Herceg" at 36 len 16)
# Para 4 "Lara Ćorić"                start=53 len=10
# Para 5 "Ana Vuksanović"            start=64 len=14
# Para 6 "Katarina Klaric"           start=79 len=15  ("Klaric" at 88 len 6)
# Para 7 "Noa Milin"                 start=95 len=9
#
# Edits are applied back-to-front (paragraph 7 first) so that earlier,
# still-to-be-processed offsets are never shifted by text inserted later
# in the run. The two in-place spelling fixes are length-preserving and
# are safe to apply at any point.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- Paragraph 7: "Noa Milin" -> append " - Dokumentacija" ---------------
$tr.Characters(95, 9).InsertAfter(" - ")
$tr.Characters(104, 3).InsertAfter("Dokumentacija")

# --- Paragraph 6: "Katarina Klaric" -> fix name, append " -  Dokumentacija"
$tr.Characters(88, 6).Text = "Klarić"
$tr.Characters(79, 15).InsertAfter(" -  ")
$tr.Characters(94, 4).InsertAfter("Dokumentacija")

# --- Paragraph 5: "Ana Vuksanović" -> append " – FE Lead " ---------------
$tr.Characters(64, 14).InsertAfter(" – FE Lead ")

# --- Paragraph 4: "Lara Ćorić" -> append " – FE Inžinjer" -----------------
$tr.Characters(53, 10).InsertAfter(" – FE ")
$tr.Characters(63, 6).InsertAfter("Inžinjer")

# --- Paragraph 3: "Mihael Breznicki-Herceg" -> fix name, append roles ----
$tr.Characters(36, 16).Text = "Breznički-Herceg"
$tr.Characters(29, 23).InsertAfter(" - ")
$tr.Characters(52, 3).InsertAfter("Dokumentacija")

# --- Paragraph 2: "Marko Pavić" -> "Marko Pavić – BE Inžinjer" -----------
$tr.Characters(17, 11).Text = "Marko Pavić – BE "
$tr.Characters(17, 17).InsertAfter("Inžinjer")

# --- Paragraph 1: "Filip Ljubotina" -> append " - Voditelj" --------------
$tr.Characters(1, 15).InsertAfter(" - ")
$tr.Characters(16, 3).InsertAfter("Voditelj")
